# Update the "想去人数" (number of people who want to go) column F values
# on both the "展览" and "全部类型" worksheets, which contain identical data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 6419
    5  = 23
    7  = 1922
    8  = 1462
    10 = 984
    11 = 305
    12 = 5600
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
